$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 2486.2144
$ws.Range("I135").Value = 2663.353
$ws.Range("J135").Value = 2212.4546
$ws.Range("K135").Value = 23970.177
$ws.Range("L135").Value = 19912.0914
$ws.Range("M135").Value = -21435.177
$ws.Range("N135").Value = -24982.0914

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 23141.068
$ws.Range("I137").Value = 26169.691
$ws.Range("J137").Value = 20680.312
$ws.Range("K137").Value = 78509.073
$ws.Range("L137").Value = 62040.936
$ws.Range("M137").Value = -75959.073
$ws.Range("N137").Value = -67140.936

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 32929.11
$ws.Range("I138").Value = 2994.1667
$ws.Range("J138").Value = 62864.055
$ws.Range("K138").Value = 8982.500100000001
$ws.Range("L138").Value = 188592.165
$ws.Range("M138").Value = -3842.500100000001
$ws.Range("N138").Value = -198872.165

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14976.963
$ws.Range("I32").Value = 15008.316
$ws.Range("K32").Value = 15008.316
$ws.Range("M32").Value = -14721.316

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 16445.215
$ws.Range("I61").Value = 8917
$ws.Range("K61").Value = 8917
$ws.Range("M61").Value = -8705

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 137651.06
$ws.Range("I74").Value = 201139.27
$ws.Range("J74").Value = 10674.667
$ws.Range("K74").Value = 201139.27
$ws.Range("L74").Value = 10674.667
$ws.Range("M74").Value = -200265.27
$ws.Range("N74").Value = -12422.667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 137651.06
$ws.Range("I77").Value = 201139.27
$ws.Range("J77").Value = 10674.667
$ws.Range("K77").Value = 1005696.35
$ws.Range("L77").Value = 53373.335
$ws.Range("M77").Value = -1001328.35
$ws.Range("N77").Value = -62109.335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2434.3513
$ws.Range("I132").Value = 2252.2646
$ws.Range("K132").Value = 6756.793799999999
$ws.Range("M132").Value = -4226.793799999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 16445.215
$ws.Range("I136").Value = 8917
$ws.Range("K136").Value = 26751
$ws.Range("M136").Value = -24201

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H138").Value = 150000
$ws.Range("J138").Value = 150000
$ws.Range("L138").Value = 150000
$ws.Range("N138").Value = -160280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3593
$ws.Range("I134").Value = 3566.1667
$ws.Range("J134").Value = 3786.2
$ws.Range("K134").Value = 10698.5001
$ws.Range("L134").Value = 11358.6
$ws.Range("M134").Value = -8163.500100000001
$ws.Range("N134").Value = -16428.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1164636.9
$ws.Range("I31").Value = 1888507.9
$ws.Range("J31").Value = 2056.182
$ws.Range("K31").Value = 1888507.9
$ws.Range("L31").Value = 2056.182
$ws.Range("M31").Value = -1888212.9
$ws.Range("N31").Value = -2646.182

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1164636.9
$ws.Range("I34").Value = 1888507.9
$ws.Range("J34").Value = 2056.182
$ws.Range("K34").Value = 1888507.9
$ws.Range("L34").Value = 2056.182
$ws.Range("M34").Value = -1888305.9
$ws.Range("N34").Value = -2460.182

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1213.5312
$ws.Range("I58").Value = 994.25
$ws.Range("K58").Value = 994.25
$ws.Range("M58").Value = -791.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 13491.647
$ws.Range("I86").Value = 13541.223
$ws.Range("K86").Value = 13541.223
$ws.Range("M86").Value = -12418.223

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 13491.647
$ws.Range("I89").Value = 13541.223
$ws.Range("K89").Value = 67706.11500000001
$ws.Range("M89").Value = -62090.11500000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 56742.777
$ws.Range("I132").Value = 91593.55
$ws.Range("K132").Value = 274780.65
$ws.Range("M132").Value = -272250.65

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1471.1777
$ws.Range("I134").Value = 1306.475
$ws.Range("K134").Value = 3919.425
$ws.Range("M134").Value = -1384.425

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1213.5312
$ws.Range("I136").Value = 994.25
$ws.Range("K136").Value = 2982.75
$ws.Range("M136").Value = -432.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 570.7727
$ws.Range("I5").Value = 528.4545000000001
$ws.Range("J5").Value = 613.0909
$ws.Range("K5").Value = 1585.3635
$ws.Range("L5").Value = 1839.2727
$ws.Range("M5").Value = -1473.3635
$ws.Range("N5").Value = -2063.2727

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3938.843
$ws.Range("J68").Value = 4416.951
$ws.Range("L68").Value = 13250.853
$ws.Range("N68").Value = -14872.853

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 3938.843
$ws.Range("J71").Value = 4416.951
$ws.Range("L71").Value = 39752.559
$ws.Range("N71").Value = -47864.559

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 11206.689
$ws.Range("I87").Value = 8748.5
$ws.Range("K87").Value = 26245.5
$ws.Range("M87").Value = -24997.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H90").Value = 11206.689
$ws.Range("I90").Value = 8748.5
$ws.Range("K90").Value = 78736.5
$ws.Range("M90").Value = -72496.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H108").Value = 1277.9231
$ws.Range("I108").Value = 984.4167
$ws.Range("J108").Value = 4800
$ws.Range("K108").Value = 2953.2501
$ws.Range("L108").Value = 14400
$ws.Range("M108").Value = -73.2501000000002
$ws.Range("N108").Value = -20160

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 570.7727
$ws.Range("I135").Value = 528.4545000000001
$ws.Range("J135").Value = 613.0909
$ws.Range("K135").Value = 4756.0905
$ws.Range("L135").Value = 5517.8181
$ws.Range("M135").Value = -2221.0905
$ws.Range("N135").Value = -10587.8181

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6284.3335
$ws.Range("I70").Value = 7640.9287
$ws.Range("J70").Value = 4823.385
$ws.Range("K70").Value = 7640.9287
$ws.Range("L70").Value = 4823.385
$ws.Range("M70").Value = -7370.9287
$ws.Range("N70").Value = -5363.385

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 6284.3335
$ws.Range("I73").Value = 7640.9287
$ws.Range("J73").Value = 4823.385
$ws.Range("K73").Value = 7640.9287
$ws.Range("L73").Value = 4823.385
$ws.Range("M73").Value = -6704.9287
$ws.Range("N73").Value = -6695.385

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1331.5555
$ws.Range("I113").Value = 1282
$ws.Range("J113").Value = 1430.6666
$ws.Range("K113").Value = 1282
$ws.Range("L113").Value = 1430.6666
$ws.Range("M113").Value = 888
$ws.Range("N113").Value = -5770.6666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 15627996
$ws.Range("I122").Value = 2897.48
$ws.Range("J122").Value = 71431920
$ws.Range("K122").Value = 8692.440000000001
$ws.Range("L122").Value = 214295760
$ws.Range("M122").Value = -6242.440000000001
$ws.Range("N122").Value = -214300660

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2362.5
$ws.Range("I126").Value = 2120
$ws.Range("K126").Value = 6360
$ws.Range("M126").Value = -3890

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2384.853
$ws.Range("I132").Value = 2200.1155
$ws.Range("K132").Value = 6600.3465
$ws.Range("M132").Value = -4070.3465

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 98995
$ws.Range("J136").Value = 98995
$ws.Range("L136").Value = 296985
$ws.Range("N136").Value = -302085

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H139").Value = 109933.336
$ws.Range("J139").Value = 109933.336
$ws.Range("L139").Value = 109933.336
$ws.Range("N139").Value = -120213.336

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4435
$ws.Range("I40").Value = 4174.1665
$ws.Range("K40").Value = 4174.1665
$ws.Range("M40").Value = -4038.1665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3238.0476
$ws.Range("I132").Value = 2727.0667
$ws.Range("K132").Value = 8181.2001
$ws.Range("M132").Value = -5651.2001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1046.5416
$ws.Range("I107").Value = 1004.1111
$ws.Range("J107").Value = 1173.8334
$ws.Range("K107").Value = 3012.3333
$ws.Range("L107").Value = 3521.5002
$ws.Range("M107").Value = -1092.3333
$ws.Range("N107").Value = -7361.5002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 25104.465
$ws.Range("I136").Value = 26971.346
$ws.Range("J136").Value = 835
$ws.Range("K136").Value = 80914.038
$ws.Range("L136").Value = 2505
$ws.Range("M136").Value = -78364.038
$ws.Range("N136").Value = -7605
